{"js": "// Apply text replacements described by the diff.\n// Each entry is [searchText, replacementText]; all source texts are\n// unique within the document, so exact, case-sensitive matches are safe.\nconst replacements = [\n  [\"2025-03-18 Tuesday\", \"2025-03-19 Wednesday\"],\n  [\"581\u00f75=\", \"176\u00f79=\"],\n  [\"704\u00f75=\", \"745\u00f75=\"],\n  [\"439\u00f76=\", \"387\u00f73=\"],\n  [\"140\u00f78=\", \"927\u00f72=\"],\n  [\"383\u00f74=\", \"546\u00f73=\"],\n  [\"964\u00f73=\", \"921\u00f72=\"],\n  [\"731\u00f73=\", \"572\u00f78=\"],\n  [\"855\u00f73=\", \"452\u00f74=\"],\n  [\"693\u00f77=\", \"962\u00f78=\"],\n  [\"208\u00f73=\", \"372\u00f72=\"],\n  [\"734\u00f76=\", \"426\u00f74=\"],\n  [\"534\u00f75=\", \"572\u00f73=\"],\n  [\"948\u00f75=\", \"403\u00f72=\"],\n  [\"436\u00f73=\", \"473\u00f75=\"],\n  [\"670\u00f77=\", \"145\u00f75=\"],\n  [\"909\u00f72=\", \"569\u00f76=\"],\n  [\"782\u00f73=\", \"418\u00f74=\"],\n  [\"553\u00f79=\", \"743\u00f76=\"],\n  [\"991\u00f76=\", \"385\u00f73=\"],\n  [\"415\u00f78=\", \"511\u00f75=\"],\n  [\"434\u00f74=\", \"687\u00f78=\"],\n  [\"529\u00f79=\", \"746\u00f77=\"],\n  [\"564\u00f72=\", \"277\u00f78=\"],\n  [\"323\u00f76=\", \"953\u00f72=\"],\n  [\"880\u00f79=\", \"414\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replacementText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply text replacements described by the diff using Word's Find/Replace.\n# Each entry is (search, replace); all source texts are unique in the\n# document, so a simple Find.Execute Replace:=wdReplaceAll (2) per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-18 Tuesday\", \"2025-03-19 Wednesday\"),\n    @(\"581\u00f75=\", \"176\u00f79=\"),\n    @(\"704\u00f75=\", \"745\u00f75=\"),\n    @(\"439\u00f76=\", \"387\u00f73=\"),\n    @(\"140\u00f78=\", \"927\u00f72=\"),\n    @(\"383\u00f74=\", \"546\u00f73=\"),\n    @(\"964\u00f73=\", \"921\u00f72=\"),\n    @(\"731\u00f73=\", \"572\u00f78=\"),\n    @(\"855\u00f73=\", \"452\u00f74=\"),\n    @(\"693\u00f77=\", \"962\u00f78=\"),\n    @(\"208\u00f73=\", \"372\u00f72=\"),\n    @(\"734\u00f76=\", \"426\u00f74=\"),\n    @(\"534\u00f75=\", \"572\u00f73=\"),\n    @(\"948\u00f75=\", \"403\u00f72=\"),\n    @(\"436\u00f73=\", \"473\u00f75=\"),\n    @(\"670\u00f77=\", \"145\u00f75=\"),\n    @(\"909\u00f72=\", \"569\u00f76=\"),\n    @(\"782\u00f73=\", \"418\u00f74=\"),\n    @(\"553\u00f79=\", \"743\u00f76=\"),\n    @(\"991\u00f76=\", \"385\u00f73=\"),\n    @(\"415\u00f78=\", \"511\u00f75=\"),\n    @(\"434\u00f74=\", \"687\u00f78=\"),\n    @(\"529\u00f79=\", \"746\u00f77=\"),\n    @(\"564\u00f72=\", \"277\u00f78=\"),\n    @(\"323\u00f76=\", \"953\u00f72=\"),\n    @(\"880\u00f79=\", \"414\u00f76=\"),\n)\n\nforeach ($pair in $replacements) {\n    $searchText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Text not found: $searchText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
